# "Generate Report for Handoff" -- refresh the localization-status report:
# a new handoff GUID (78ba8c9a-a0f4-43fe-a795-cfdec66eaf74) replaces the
# previous one (9851bec4-9180-445f-b2aa-09b1645016d2), and the handoff /
# xliff-generate timestamps move forward a few seconds.

$wb = $excel.ActiveWorkbook

$oldGuid = "9851bec4-9180-445f-b2aa-09b1645016d2"
$newGuid = "78ba8c9a-a0f4-43fe-a795-cfdec66eaf74"

# The external hyperlink targets (same repo blob URL) are unchanged by this
# edit -- only the visible display text / cell text move to the new GUID.
$linkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0080857e036589d4ba1b69b3d2dfed5e27235f15/e2e/$oldGuid.md"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $linkTarget, "", "", "e2e\$newGuid.md")

$wsOverview.Range("G2").Value = "2016-09-06 15:57:20"

$wsOverview.Columns.Item(1).ColumnWidth = 39

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $linkTarget, "", "", "$newGuid.md")

$wsZhCn.Range("G2").Value = "$newGuid.6f012e761a17d4d2eb37dff6a0e3b8433ed462c1.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-06 15:57:02"

$wsZhCn.Columns.Item(1).ColumnWidth = 39

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $linkTarget, "", "", "$newGuid.md")

$wsDeDe.Range("G2").Value = "$newGuid.6f012e761a17d4d2eb37dff6a0e3b8433ed462c1.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-06 15:57:20"

$wsDeDe.Columns.Item(1).ColumnWidth = 39
